$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.367.97'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '1.874.39'
$ws.Range("E3").Value = '  +0.65%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.0000'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7123'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.97'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.65%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07803'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.57%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3116'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.18'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08448'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.81%  '
$ws.Range("D12").Value = '1.872.94'
$ws.Range("E12").Value = '  +0.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.237'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7134'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.10'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.06%  '
$ws.Range("D16").Value = '29.374.25'
$ws.Range("E16").Value = '  +0.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.065'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.92%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008231'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '240.89'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.58%  '
$ws.Range("E20").Value = '  +0.71%  '
$ws.Range("D21").Value = '2.122.15'
$ws.Range("E21").Value = '  +0.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.791'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.50%  '
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1594'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.069'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.56'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.512'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.422'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.43%  '
$ws.Range("E31").Value = '  -3.68%  '
$ws.Range("E32").Value = '  +2.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05309'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.939'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.50%  '
$ws.Range("E35").Value = '  +1.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7450'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.93%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.697'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01868'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.02%  '
$ws.Range("D39").Value = '1.226.19'
$ws.Range("E39").Value = '  +4.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.504'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '110.73'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +8.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8910'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.20%  '
$ws.Range("E44").Value = '  +0.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.0000'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").Value = '2.019.93'
$ws.Range("E46").Value = '  +0.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.812'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.59%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5213'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.65%  '
$ws.Range("E49").Value = '  +2.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.440'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4324'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.41%  '
